# Scheduled-runner refresh: updates computed market/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) for a set of leve rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 168624.83
$ws.Range("I111").Value = 1649
$ws.Range("J111").Value = 202020
$ws.Range("K111").Value = 4947
$ws.Range("L111").Value = 606060
$ws.Range("M111").Value = -1880
$ws.Range("N111").Value = -612194
$ws.Range("H113").Value = 6767.5
$ws.Range("I113").Value = 4276.25
$ws.Range("J113").Value = 11750
$ws.Range("K113").Value = 4276.25
$ws.Range("L113").Value = 11750
$ws.Range("M113").Value = -1022.25
$ws.Range("N113").Value = -18258

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4224.5713
$ws.Range("I2").Value = 4914.4
$ws.Range("J2").Value = 2500
$ws.Range("K2").Value = 4914.4
$ws.Range("L2").Value = 2500
$ws.Range("M2").Value = -4801.4
$ws.Range("N2").Value = -2726
$ws.Range("H45").Value = 9343.6
$ws.Range("I45").Value = 17368.416
$ws.Range("J45").Value = 1936.0769
$ws.Range("K45").Value = 17368.416
$ws.Range("L45").Value = 1936.0769
$ws.Range("M45").Value = -16991.416
$ws.Range("N45").Value = -2690.0769
$ws.Range("H97").Value = 744.381
$ws.Range("I97").Value = 542.3077
$ws.Range("J97").Value = 1072.75
$ws.Range("K97").Value = 542.3077
$ws.Range("L97").Value = 1072.75
$ws.Range("M97").Value = -46.30769999999995
$ws.Range("N97").Value = -2064.75
$ws.Range("H102").Value = 9261959
$ws.Range("I102").Value = 12347612
$ws.Range("K102").Value = 12347612
$ws.Range("M102").Value = -12345990
$ws.Range("H110").Value = 876.8461
$ws.Range("I110").Value = 857
$ws.Range("K110").Value = 857
$ws.Range("M110").Value = 1188
$ws.Range("H116").Value = 4224.5713
$ws.Range("I116").Value = 4914.4
$ws.Range("J116").Value = 2500
$ws.Range("K116").Value = 4914.4
$ws.Range("L116").Value = 2500
$ws.Range("M116").Value = -2620.4
$ws.Range("N116").Value = -7088
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").ClearContents()
$ws.Range("N117").Value = 0
$ws.Range("H122").Value = 6412209
$ws.Range("I122").Value = 8548612
$ws.Range("K122").Value = 25645836
$ws.Range("M122").Value = -25643386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4224.5713
$ws.Range("I3").Value = 4914.4
$ws.Range("J3").Value = 2500
$ws.Range("K3").Value = 4914.4
$ws.Range("L3").Value = 2500
$ws.Range("M3").Value = -4800.4
$ws.Range("N3").Value = -2728
$ws.Range("H94").Value = 2072.0454
$ws.Range("I94").Value = 1862.8572
$ws.Range("J94").Value = 2169.6667
$ws.Range("K94").Value = 1862.8572
$ws.Range("L94").Value = 2169.6667
$ws.Range("M94").Value = -1411.8572
$ws.Range("N94").Value = -3071.6667
$ws.Range("H99").Value = 333334800
$ws.Range("I99").Value = 500000960
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 500000960
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = -499999462
$ws.Range("N99").Value = -5496
$ws.Range("H100").Value = 41000
$ws.Range("J100").Value = 41000
$ws.Range("L100").Value = 41000
$ws.Range("N100").Value = -43164
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("N103").Value = 0
$ws.Range("H107").Value = 1199
$ws.Range("I107").Value = 1196.1
$ws.Range("J107").Value = 1208.6666
$ws.Range("K107").Value = 1196.1
$ws.Range("L107").Value = 1208.6666
$ws.Range("M107").Value = 723.9000000000001
$ws.Range("N107").Value = -5048.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 39603.383
$ws.Range("J106").Value = 39603.383
$ws.Range("L106").Value = 39603.383
$ws.Range("N106").Value = -42127.383
$ws.Range("H122").Value = 1042.6154
$ws.Range("I122").Value = 1213.5
$ws.Range("J122").Value = 966.6667
$ws.Range("K122").Value = 3640.5
$ws.Range("L122").Value = 2900.0001
$ws.Range("M122").Value = -1190.5
$ws.Range("N122").Value = -7800.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 429765
$ws.Range("I5").Value = 786
$ws.Range("J5").Value = 668086.7
$ws.Range("K5").Value = 2358
$ws.Range("L5").Value = 2004260.1
$ws.Range("M5").Value = -2246
$ws.Range("N5").Value = -2004484.1
$ws.Range("H122").Value = 977.6667
$ws.Range("I122").Value = 800
$ws.Range("J122").Value = 1199.75
$ws.Range("K122").Value = 7200
$ws.Range("L122").Value = 10797.75
$ws.Range("M122").Value = -4750
$ws.Range("N122").Value = -15697.75
$ws.Range("H135").Value = 429765
$ws.Range("I135").Value = 786
$ws.Range("J135").Value = 668086.7
$ws.Range("K135").Value = 7074
$ws.Range("L135").Value = 6012780.3
$ws.Range("M135").Value = -4539
$ws.Range("N135").Value = -6017850.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2161.2
$ws.Range("I97").Value = 2161.2
$ws.Range("K97").Value = 2161.2
$ws.Range("M97").Value = -1665.2
$ws.Range("H113").Value = 200000850
$ws.Range("I113").Value = 200000850
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 200000850
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -199998680
$ws.Range("H122").Value = 12758661
$ws.Range("I122").Value = 10803875
$ws.Range("J122").Value = 16668233
$ws.Range("K122").Value = 32411625
$ws.Range("L122").Value = 50004699
$ws.Range("M122").Value = -32409175
$ws.Range("N122").Value = -50009599

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2250.3333
$ws.Range("I61").Value = 2289.3333
$ws.Range("J61").Value = 2133.3333
$ws.Range("K61").Value = 2289.3333
$ws.Range("L61").Value = 2133.3333
$ws.Range("M61").Value = -2087.3333
$ws.Range("N61").Value = -2537.3333
$ws.Range("H100").Value = 3365.4
$ws.Range("I100").Value = 3001.5
$ws.Range("K100").Value = 3001.5
$ws.Range("M100").Value = -2460.5
$ws.Range("H113").Value = 2250.3333
$ws.Range("I113").Value = 2289.3333
$ws.Range("J113").Value = 2133.3333
$ws.Range("K113").Value = 2289.3333
$ws.Range("L113").Value = 2133.3333
$ws.Range("M113").Value = -119.3332999999998
$ws.Range("N113").Value = -6473.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 52632068
$ws.Range("I107").Value = 76923580
$ws.Range("J107").Value = 446.66666
$ws.Range("K107").Value = 230770740
$ws.Range("L107").Value = 1339.99998
$ws.Range("M107").Value = -230768820
$ws.Range("N107").Value = -5179.999980000001
